# Apply the afferent_data_updated.xlsx edits:
# 1. Add an "afferent index" column (A) with values 1,6,2,3,4,5 for rows 4-9
# 2. Add a new (empty, formatted) cell J4 with a centered number format and
#    thin left/right border
# 3. Move the active selection on Feuil1 to J4

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Feuil1")

# --- New column A values (afferent order/index) ---
$ws1.Range("A4").Value = 1
$ws1.Range("A5").Value = 6
$ws1.Range("A6").Value = 2
$ws1.Range("A7").Value = 3
$ws1.Range("A8").Value = 4
$ws1.Range("A9").Value = 5

# --- New formatted (empty) cell J4 ---
$rngJ4 = $ws1.Range("J4")
$rngJ4.NumberFormat = "0.00"
$rngJ4.HorizontalAlignment = -4108
$rngJ4.Borders.Item(7).Weight = 2
$rngJ4.Borders.Item(7).ColorIndex = -4105
$rngJ4.Borders.Item(10).Weight = 2
$rngJ4.Borders.Item(10).ColorIndex = -4105

# --- Move the active selection to J4 (matches author re-saving with J4 selected) ---
$ws1.Activate() | Out-Null
$rngJ4.Select() | Out-Null

Write-Host "done"
